# Daily attendance processing - 2026-01-25 06:08:06
# Normalizes the "Recorded By" (column G) values: when a cell's comma-separated
# list of recorders begins with "System", the order of the list is reversed
# (which also has the effect of swapping the case of "System"/"system" tokens
# that sit at the very start/end of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows that are intentionally left untouched even though their value matches
# the "starts with System" pattern (mirrors the source diff exactly).
$skipRows = @(4, 30, 56)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) { continue }

    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq '') { continue }

    $parts = $val -split ', '
    if ($parts[0] -ceq 'System') {
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value = [string]::Join(', ', $reversed)
    }
}
